$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Iwan Setiyawan's row (row 6) Nama/TTL text: drop the Bantul/2022-06-06 suffix
$ws.Range("B6").Value = "Iwan Setiyawan//"

# Clear "Pendidikan Terakhir" (education level) values that are being removed
# from the workbook (S1 / S2 / SMP / SMK / DIII / SMA all go away).
$eRowsToClear = @(4,5,6,7,8,10,11,12,13,15,16,17,18,20,31,148)
foreach ($r in $eRowsToClear) {
    $ws.Cells.Item($r, 5).ClearContents()
}

# Clear "Jabatan" values that referenced removed shared strings
# (Administrasi Personalia on row 148, Staff on row 168).
$ws.Cells.Item(148, 4).ClearContents()
$ws.Cells.Item(168, 4).ClearContents()

# Row 168 (B/ratnasari's NIP-like numeric column) changes value
$ws.Range("C168").Value = 9535253

# Row 169 (ratnasari/Darit/2001-08-09) is removed entirely
$ws.Rows.Item(169).Delete()
